# Fix: prevent rows whose "before" (FV2210) and "after" (FV2304) halves are
# identical from still being flagged as "ÄNDERUNG" (changed) in column L.
# Also restores the group-header highlighting (gray fill, bold segment name)
# on the first row of each newly-unflagged group, matching the styling that
# correctly-handled rows elsewhere in the sheet already use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastRow = $ws.UsedRange.Rows.Count

$leftCols  = @("B","C","D","E","F","G","H","I","J","K")
$rightCols = @("M","N","O","P","Q","R","S","T","U","V")

$GRAY = 14277081   # fill RGB D9D9D9 (BGR-packed long used by COM Interior.Color)

# Track the previous data row's "segment name" (column B) so we can tell
# whether a row starts a new group (used to decide whether to re-apply the
# bold/gray "group header" look once its bogus change-flag is cleared).
$prevB = $null

for ($r = $firstDataRow; $r -le $lastRow; $r++) {

    $same = $true
    for ($i = 0; $i -lt $leftCols.Count; $i++) {
        $v1 = $ws.Range($leftCols[$i] + $r).Value2
        $v2 = $ws.Range($rightCols[$i] + $r).Value2
        if ($v1 -eq $null -and $v2 -eq $null) {
            continue
        }
        if ($v1 -ne $v2) {
            $same = $false
            break
        }
    }

    $bVal = $ws.Range("B" + $r).Value2
    $isGroupStart = ($bVal -ne $prevB)

    if ($same) {
        $lCell = $ws.Range("L" + $r)
        if ($lCell.Value2 -ne $null) {
            # This row has no real difference between the two AHB versions,
            # so the "ÄNDERUNG" marker is a false positive - remove it.
            $lCell.ClearContents()
        }
        # Column L keeps the "empty" look: gray fill, centered, regular (not
        # bold/colored) text - same formatting already used by every other
        # un-flagged row in the sheet.
        $lCell.Interior.Color = $GRAY
        $lCell.Font.Bold = $false
        $lCell.Font.Color = 0
        $lCell.HorizontalAlignment = -4108  # xlCenter
        $lCell.WrapText = $true

        if ($isGroupStart) {
            # First row of a segment/segment-group: highlight the whole row
            # (gray fill across A:V, bold segment name in column B) like the
            # other correctly-rendered group headers in the sheet.
            $rowRange = $ws.Range("A" + $r + ":V" + $r)
            $rowRange.Interior.Color = $GRAY
            $rowRange.Font.Bold = $false
            $rowRange.WrapText = $true

            $bCell = $ws.Range("B" + $r)
            $bCell.Font.Bold = $true
        }
    }

    $prevB = $bVal
}
